$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the formatting of the existing
# header cells (bold, centered, bordered) by copying G1's format.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the corresponding data value for row 2 (Save = 0)
$ws.Range("H2").Value = 0
